$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Gesellschaft" -> "Unternehmen" and "Gesellschaftskuerzel" -> "Unternehmenskuerzel"
$ws.Range("A2").Value = "Unternehmen"
$ws.Range("A3").Value = "Unternehmenskuerzel"

# Update the active selection to A4 (was B3)
$ws.Range("A4").Select()
